$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 2290.6155  # H100: 2513.9092 -> 2290.6155
$ws.Cells.Item(100, 9).Value = 2397.25  # I100: 2842.1667 -> 2397.25
$ws.Cells.Item(100, 11).Value = 2397.25  # K100: 2842.1667 -> 2397.25
$ws.Cells.Item(100, 13).Value = -1856.25  # M100: -2301.1667 -> -1856.25

$ws.Cells.Item(112, 8).Value = 2895.3  # H112: 2955.3 -> 2895.3
$ws.Cells.Item(112, 10).Value = 2681.625  # J112: 2756.625 -> 2681.625
$ws.Cells.Item(112, 12).Value = 8044.875  # L112: 8269.875 -> 8044.875
$ws.Cells.Item(112, 14).Value = -10260.875  # N112: -10485.875 -> -10260.875

$ws.Cells.Item(135, 8).Value = 1260.6923  # H135: 1314.5 -> 1260.6923
$ws.Cells.Item(135, 9).Value = 1123  # I135: 1195.3 -> 1123
$ws.Cells.Item(135, 10).Value = 2018  # J135: 1910.5 -> 2018
$ws.Cells.Item(135, 11).Value = 10107  # K135: 10757.7 -> 10107
$ws.Cells.Item(135, 12).Value = 18162  # L135: 17194.5 -> 18162
$ws.Cells.Item(135, 13).Value = -7572  # M135: -8222.699999999999 -> -7572
$ws.Cells.Item(135, 14).Value = -23232  # N135: -22264.5 -> -23232

$ws.Cells.Item(137, 8).Value = 1571.6  # H137: 1635.1666 -> 1571.6
$ws.Cells.Item(137, 9).Value = 1168.1428  # I137: 1196.25 -> 1168.1428
$ws.Cells.Item(137, 11).Value = 3504.4284  # K137: 3588.75 -> 3504.4284
$ws.Cells.Item(137, 13).Value = -954.4284000000002  # M137: -1038.75 -> -954.4284000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7694744  # H32: 9093630 -> 7694744
$ws.Cells.Item(32, 9).Value = 2639.5833  # I32: 2992.5 -> 2639.5833
$ws.Cells.Item(32, 11).Value = 2639.5833  # K32: 2992.5 -> 2639.5833
$ws.Cells.Item(32, 13).Value = -2352.5833  # M32: -2705.5 -> -2352.5833

$ws.Cells.Item(54, 8).Value = 16000  # H54: 0 -> 16000
$ws.Cells.Item(54, 10).Value = 16000  # J54: 0 -> 16000
$ws.Cells.Item(54, 12).Value = 16000  # L54: 0 -> 16000
$ws.Cells.Item(54, 14).Value = -17538  # N54: None -> -17538

$ws.Cells.Item(74, 8).Value = 7540.2856  # H74: 7360.5 -> 7540.2856
$ws.Cells.Item(74, 9).Value = 7273.25  # I74: 7130.857 -> 7273.25
$ws.Cells.Item(74, 11).Value = 7273.25  # K74: 7130.857 -> 7273.25
$ws.Cells.Item(74, 13).Value = -6399.25  # M74: -6256.857 -> -6399.25

$ws.Cells.Item(77, 8).Value = 7540.2856  # H77: 7360.5 -> 7540.2856
$ws.Cells.Item(77, 9).Value = 7273.25  # I77: 7130.857 -> 7273.25
$ws.Cells.Item(77, 11).Value = 36366.25  # K77: 35654.285 -> 36366.25
$ws.Cells.Item(77, 13).Value = -31998.25  # M77: -31286.285 -> -31998.25

$ws.Cells.Item(92, 8).Value = 47000  # H92: 0 -> 47000
$ws.Cells.Item(92, 10).Value = 47000  # J92: 0 -> 47000
$ws.Cells.Item(92, 12).Value = 47000  # L92: 0 -> 47000
$ws.Cells.Item(92, 14).Value = -51992  # N92: None -> -51992

$ws.Cells.Item(96, 8).Value = 2243382  # H96: 2880541 -> 2243382
$ws.Cells.Item(96, 10).Value = 2243382  # J96: 2880541 -> 2243382
$ws.Cells.Item(96, 12).Value = 2243382  # L96: 2880541 -> 2243382
$ws.Cells.Item(96, 14).Value = -2248874  # N96: -2886033 -> -2248874

$ws.Cells.Item(97, 8).Value = 704.0833  # H97: 777.5 -> 704.0833
$ws.Cells.Item(97, 10).Value = 337  # J97: 0 -> 337
$ws.Cells.Item(97, 12).Value = 337  # L97: 0 -> 337
$ws.Cells.Item(97, 14).Value = -1329  # N97: None -> -1329

$ws.Cells.Item(104, 8).Value = 15658.333  # H104: 21987.5 -> 15658.333
$ws.Cells.Item(104, 10).Value = 15658.333  # J104: 21987.5 -> 15658.333
$ws.Cells.Item(104, 12).Value = 15658.333  # L104: 21987.5 -> 15658.333
$ws.Cells.Item(104, 14).Value = -22646.333  # N104: -28975.5 -> -22646.333

$ws.Cells.Item(122, 8).Value = 1360.875  # H122: 1374.25 -> 1360.875
$ws.Cells.Item(122, 9).Value = 1126.8572  # I122: 1149.1666 -> 1126.8572
$ws.Cells.Item(122, 10).Value = 2999  # J122: 2049.5 -> 2999
$ws.Cells.Item(122, 11).Value = 3380.5716  # K122: 3447.4998 -> 3380.5716
$ws.Cells.Item(122, 12).Value = 8997  # L122: 6148.5 -> 8997
$ws.Cells.Item(122, 13).Value = -930.5715999999998  # M122: -997.4998000000001 -> -930.5715999999998
$ws.Cells.Item(122, 14).Value = -13897  # N122: -11048.5 -> -13897

$ws.Cells.Item(132, 8).Value = 1746.5  # H132: 1639.8572 -> 1746.5
$ws.Cells.Item(132, 9).Value = 1695.8  # I132: 1579.8334 -> 1695.8
$ws.Cells.Item(132, 11).Value = 5087.4  # K132: 4739.5002 -> 5087.4
$ws.Cells.Item(132, 13).Value = -2557.4  # M132: -2209.5002 -> -2557.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 1389135.8  # H22: 1234805.1 -> 1389135.8
$ws.Cells.Item(22, 9).Value = 3703817  # I22: 2222349.2 -> 3703817
$ws.Cells.Item(22, 10).Value = 327  # J22: 375 -> 327
$ws.Cells.Item(22, 11).Value = 3703817  # K22: 2222349.2 -> 3703817
$ws.Cells.Item(22, 12).Value = 327  # L22: 375 -> 327
$ws.Cells.Item(22, 13).Value = -3703644  # M22: -2222176.2 -> -3703644
$ws.Cells.Item(22, 14).Value = -673  # N22: -721 -> -673

$ws.Cells.Item(24, 8).Value = 900  # H24: 1166.6666 -> 900
$ws.Cells.Item(24, 9).Value = 800  # I24: 1166.6666 -> 800
$ws.Cells.Item(24, 10).Value = 1000  # J24: 0 -> 1000
$ws.Cells.Item(24, 11).Value = 800  # K24: 1166.6666 -> 800
$ws.Cells.Item(24, 12).Value = 1000  # L24: 0 -> 1000
$ws.Cells.Item(24, 13).Value = -565  # M24: -931.6666 -> -565
$ws.Cells.Item(24, 14).Value = -1470  # N24: None -> -1470

$ws.Cells.Item(36, 8).Value = 887.125  # H36: 988.8570999999999 -> 887.125
$ws.Cells.Item(36, 9).Value = 887.125  # I36: 988.8570999999999 -> 887.125
$ws.Cells.Item(36, 11).Value = 887.125  # K36: 988.8570999999999 -> 887.125
$ws.Cells.Item(36, 13).Value = -353.125  # M36: -454.8570999999999 -> -353.125

$ws.Cells.Item(134, 8).Value = 5822.615  # H134: 5823.077 -> 5822.615
$ws.Cells.Item(134, 9).Value = 1190.3636  # I134: 1190.909 -> 1190.3636
$ws.Cells.Item(134, 11).Value = 3571.0908  # K134: 3572.727 -> 3571.0908
$ws.Cells.Item(134, 13).Value = -1036.0908  # M134: -1037.727 -> -1036.0908

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5124.448  # H31: 4980.2666 -> 5124.448
$ws.Cells.Item(31, 9).Value = 2161.1  # I31: 2037.2727 -> 2161.1
$ws.Cells.Item(31, 11).Value = 2161.1  # K31: 2037.2727 -> 2161.1
$ws.Cells.Item(31, 13).Value = -1866.1  # M31: -1742.2727 -> -1866.1

$ws.Cells.Item(34, 8).Value = 5124.448  # H34: 4980.2666 -> 5124.448
$ws.Cells.Item(34, 9).Value = 2161.1  # I34: 2037.2727 -> 2161.1
$ws.Cells.Item(34, 11).Value = 2161.1  # K34: 2037.2727 -> 2161.1
$ws.Cells.Item(34, 13).Value = -1959.1  # M34: -1835.2727 -> -1959.1

$ws.Cells.Item(134, 8).Value = 3010.25  # H134: 3164.2666 -> 3010.25
$ws.Cells.Item(134, 9).Value = 2885.0667  # I134: 3041.1428 -> 2885.0667
$ws.Cells.Item(134, 11).Value = 8655.2001  # K134: 9123.428400000001 -> 8655.2001
$ws.Cells.Item(134, 13).Value = -6120.2001  # M134: -6588.428400000001 -> -6120.2001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1933.6428  # H5: 1959.3077 -> 1933.6428
$ws.Cells.Item(5, 10).Value = 2517.8  # J5: 2747.25 -> 2517.8
$ws.Cells.Item(5, 12).Value = 7553.400000000001  # L5: 8241.75 -> 7553.400000000001
$ws.Cells.Item(5, 14).Value = -7777.400000000001  # N5: -8465.75 -> -7777.400000000001

$ws.Cells.Item(40, 8).Value = 240.375  # H40: 214.77777 -> 240.375
$ws.Cells.Item(40, 10).Value = 600  # J40: 452.5 -> 600
$ws.Cells.Item(40, 12).Value = 2400  # L40: 1810 -> 2400
$ws.Cells.Item(40, 14).Value = -2538  # N40: -1948 -> -2538

$ws.Cells.Item(80, 8).Value = 4372.6816  # H80: 4192.8 -> 4372.6816
$ws.Cells.Item(80, 9).Value = 4066.611  # I80: 3959.95 -> 4066.611
$ws.Cells.Item(80, 10).Value = 5750  # J80: 5124.2 -> 5750
$ws.Cells.Item(80, 11).Value = 12199.833  # K80: 11879.85 -> 12199.833
$ws.Cells.Item(80, 12).Value = 17250  # L80: 15372.6 -> 17250
$ws.Cells.Item(80, 13).Value = -11263.833  # M80: -10943.85 -> -11263.833
$ws.Cells.Item(80, 14).Value = -19122  # N80: -17244.6 -> -19122

$ws.Cells.Item(83, 8).Value = 4372.6816  # H83: 4192.8 -> 4372.6816
$ws.Cells.Item(83, 9).Value = 4066.611  # I83: 3959.95 -> 4066.611
$ws.Cells.Item(83, 10).Value = 5750  # J83: 5124.2 -> 5750
$ws.Cells.Item(83, 11).Value = 36599.499  # K83: 35639.55 -> 36599.499
$ws.Cells.Item(83, 12).Value = 51750  # L83: 46117.8 -> 51750
$ws.Cells.Item(83, 13).Value = -31919.499  # M83: -30959.55 -> -31919.499
$ws.Cells.Item(83, 14).Value = -61110  # N83: -55477.8 -> -61110

$ws.Cells.Item(107, 8).Value = 452.51428  # H107: 458.45715 -> 452.51428
$ws.Cells.Item(107, 9).Value = 213.22223  # I107: 213.77777 -> 213.22223
$ws.Cells.Item(107, 10).Value = 705.8823  # J107: 717.5294 -> 705.8823
$ws.Cells.Item(107, 11).Value = 639.66669  # K107: 641.33331 -> 639.66669
$ws.Cells.Item(107, 12).Value = 2117.6469  # L107: 2152.5882 -> 2117.6469
$ws.Cells.Item(107, 13).Value = 1280.33331  # M107: 1278.66669 -> 1280.33331
$ws.Cells.Item(107, 14).Value = -5957.6469  # N107: -5992.5882 -> -5957.6469

$ws.Cells.Item(116, 8).Value = 1859  # H116: 2099.5 -> 1859
$ws.Cells.Item(116, 9).Value = 1948.75  # I116: 2300 -> 1948.75
$ws.Cells.Item(116, 10).Value = 1500  # J116: 1498 -> 1500
$ws.Cells.Item(116, 11).Value = 5846.25  # K116: 6900 -> 5846.25
$ws.Cells.Item(116, 12).Value = 4500  # L116: 4494 -> 4500
$ws.Cells.Item(116, 13).Value = -2404.25  # M116: -3458 -> -2404.25
$ws.Cells.Item(116, 14).Value = -11384  # N116: -11378 -> -11384

$ws.Cells.Item(131, 8).Value = 2042.9231  # H131: 2238.3333 -> 2042.9231
$ws.Cells.Item(131, 9).Value = 1549.25  # I131: 1832.3334 -> 1549.25
$ws.Cells.Item(131, 10).Value = 2262.3333  # J131: 2373.6667 -> 2262.3333
$ws.Cells.Item(131, 11).Value = 4647.75  # K131: 5497.0002 -> 4647.75
$ws.Cells.Item(131, 12).Value = 6786.999899999999  # L131: 7121.000100000001 -> 6786.999899999999
$ws.Cells.Item(131, 13).Value = 392.25  # M131: -457.0002000000004 -> 392.25
$ws.Cells.Item(131, 14).Value = -16866.9999  # N131: -17201.0001 -> -16866.9999

$ws.Cells.Item(135, 8).Value = 1933.6428  # H135: 1959.3077 -> 1933.6428
$ws.Cells.Item(135, 10).Value = 2517.8  # J135: 2747.25 -> 2517.8
$ws.Cells.Item(135, 12).Value = 22660.2  # L135: 24725.25 -> 22660.2
$ws.Cells.Item(135, 14).Value = -27730.2  # N135: -29795.25 -> -27730.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1859.4166  # H102: 1566.1765 -> 1859.4166
$ws.Cells.Item(102, 9).Value = 955.375  # I102: 829.5833 -> 955.375
$ws.Cells.Item(102, 10).Value = 3667.5  # J102: 3334 -> 3667.5
$ws.Cells.Item(102, 11).Value = 955.375  # K102: 829.5833 -> 955.375
$ws.Cells.Item(102, 12).Value = 3667.5  # L102: 3334 -> 3667.5
$ws.Cells.Item(102, 13).Value = 666.625  # M102: 792.4167 -> 666.625
$ws.Cells.Item(102, 14).Value = -6911.5  # N102: -6578 -> -6911.5

$ws.Cells.Item(132, 8).Value = 79181.38  # H132: 85628.75 -> 79181.38
$ws.Cells.Item(132, 9).Value = 85596.5  # I132: 102314.6 -> 85596.5
$ws.Cells.Item(132, 10).Value = 2200  # J132: 2199.5 -> 2200
$ws.Cells.Item(132, 11).Value = 256789.5  # K132: 306943.8 -> 256789.5
$ws.Cells.Item(132, 12).Value = 6600  # L132: 6598.5 -> 6600
$ws.Cells.Item(132, 13).Value = -254259.5  # M132: -304413.8 -> -254259.5
$ws.Cells.Item(132, 14).Value = -11660  # N132: -11658.5 -> -11660

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1750  # H22: 1866.6666 -> 1750
$ws.Cells.Item(22, 9).Value = 1000  # I22: 1550 -> 1000
$ws.Cells.Item(22, 11).Value = 1000  # K22: 1550 -> 1000
$ws.Cells.Item(22, 13).Value = -705  # M22: -1255 -> -705

$ws.Cells.Item(27, 8).Value = 1750  # H27: 1866.6666 -> 1750
$ws.Cells.Item(27, 9).Value = 1000  # I27: 1550 -> 1000
$ws.Cells.Item(27, 11).Value = 1000  # K27: 1550 -> 1000
$ws.Cells.Item(27, 13).Value = -893  # M27: -1443 -> -893

$ws.Cells.Item(68, 8).Value = 4518.643  # H68: 5036.4287 -> 4518.643
$ws.Cells.Item(68, 9).Value = 2108.7144  # I68: 2335.1667 -> 2108.7144
$ws.Cells.Item(68, 10).Value = 6928.5713  # J68: 7062.375 -> 6928.5713
$ws.Cells.Item(68, 11).Value = 2108.7144  # K68: 2335.1667 -> 2108.7144
$ws.Cells.Item(68, 12).Value = 6928.5713  # L68: 7062.375 -> 6928.5713
$ws.Cells.Item(68, 13).Value = -1359.7144  # M68: -1586.1667 -> -1359.7144
$ws.Cells.Item(68, 14).Value = -8426.5713  # N68: -8560.375 -> -8426.5713

$ws.Cells.Item(71, 8).Value = 4518.643  # H71: 5036.4287 -> 4518.643
$ws.Cells.Item(71, 9).Value = 2108.7144  # I71: 2335.1667 -> 2108.7144
$ws.Cells.Item(71, 10).Value = 6928.5713  # J71: 7062.375 -> 6928.5713
$ws.Cells.Item(71, 11).Value = 10543.572  # K71: 11675.8335 -> 10543.572
$ws.Cells.Item(71, 12).Value = 34642.85649999999  # L71: 35311.875 -> 34642.85649999999
$ws.Cells.Item(71, 13).Value = -6799.572  # M71: -7931.833500000001 -> -6799.572
$ws.Cells.Item(71, 14).Value = -42130.85649999999  # N71: -42799.875 -> -42130.85649999999

$ws.Cells.Item(101, 8).Value = 13751.833  # H101: 14138.363 -> 13751.833
$ws.Cells.Item(101, 10).Value = 14063  # J101: 14519.3 -> 14063
$ws.Cells.Item(101, 12).Value = 14063  # L101: 14519.3 -> 14063
$ws.Cells.Item(101, 14).Value = -20553  # N101: -21009.3 -> -20553

$ws.Cells.Item(132, 8).Value = 3289.6  # H132: 3354.182 -> 3289.6
$ws.Cells.Item(132, 10).Value = 3401.25  # J132: 3521 -> 3401.25
$ws.Cells.Item(132, 12).Value = 10203.75  # L132: 10563 -> 10203.75
$ws.Cells.Item(132, 14).Value = -15263.75  # N132: -15623 -> -15263.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(58, 8).Value = 3128.125  # H58: 3112.2222 -> 3128.125
$ws.Cells.Item(58, 9).Value = 3128.125  # I58: 3112.2222 -> 3128.125
$ws.Cells.Item(58, 11).Value = 3128.125  # K58: 3112.2222 -> 3128.125
$ws.Cells.Item(58, 13).Value = -2820.125  # M58: -2804.2222 -> -2820.125

$ws.Cells.Item(110, 8).Value = 0  # H110: 5000 -> 0
$ws.Cells.Item(110, 10).Value = 0  # J110: 5000 -> 0
$ws.Cells.Item(110, 12).Value = 0  # L110: 5000 -> 0
$ws.Cells.Item(110, 14).ClearContents()  # N110: -13180 -> (removed)

$ws.Cells.Item(132, 8).Value = 1234.1111  # H132: 1469.4667 -> 1234.1111
$ws.Cells.Item(132, 9).Value = 1027.6666  # I132: 1086.9166 -> 1027.6666
$ws.Cells.Item(132, 10).Value = 2266.3333  # J132: 2999.6667 -> 2266.3333
$ws.Cells.Item(132, 11).Value = 3082.9998  # K132: 3260.7498 -> 3082.9998
$ws.Cells.Item(132, 12).Value = 6798.999899999999  # L132: 8999.000100000001 -> 6798.999899999999
$ws.Cells.Item(132, 13).Value = -552.9998000000001  # M132: -730.7498000000001 -> -552.9998000000001
$ws.Cells.Item(132, 14).Value = -11858.9999  # N132: -14059.0001 -> -11858.9999

$ws.Cells.Item(136, 8).Value = 2887.647  # H136: 3064.1177 -> 2887.647
$ws.Cells.Item(136, 9).Value = 1836.7778  # I136: 1941.625 -> 1836.7778
$ws.Cells.Item(136, 10).Value = 4069.875  # J136: 4061.889 -> 4069.875
$ws.Cells.Item(136, 11).Value = 5510.3334  # K136: 5824.875 -> 5510.3334
$ws.Cells.Item(136, 12).Value = 12209.625  # L136: 12185.667 -> 12209.625
$ws.Cells.Item(136, 13).Value = -2960.3334  # M136: -3274.875 -> -2960.3334
$ws.Cells.Item(136, 14).Value = -17309.625  # N136: -17285.667 -> -17309.625
